# Updates the cryptos price/volume table (Coin/Link/Price/Volume(1h) columns B:E)
# to the latest scraped values. Price cells whose text looks like a plain
# decimal number are prefixed with a leading apostrophe so Excel keeps them
# as literal text (matching the workbook's original inlineStr/text cells)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.585.00'
$ws.Cells.Item(2, 5).Value = '  +1.20%  '
$ws.Cells.Item(3, 4).Value = '3.347.32'
$ws.Cells.Item(3, 5).Value = '  +1.27%  '
$ws.Cells.Item(4, 4).Value = '''0.998'
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = '''586.16'
$ws.Cells.Item(5, 5).Value = '  +5.76%  '
$ws.Cells.Item(6, 4).Value = '''186.83'
$ws.Cells.Item(6, 5).Value = '  -0.14%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '''0.577'
$ws.Cells.Item(8, 5).Value = '  -0.48%  '
$ws.Cells.Item(9, 5).Value = '  +1.47%  '
$ws.Cells.Item(10, 4).Value = '''0.584'
$ws.Cells.Item(10, 5).Value = '  +0.71%  '
$ws.Cells.Item(11, 4).Value = '''47.21'
$ws.Cells.Item(11, 5).Value = '  +0.94%  '
$ws.Cells.Item(12, 4).Value = '''0.0000273'
$ws.Cells.Item(12, 5).Value = '  +2.17%  '
$ws.Cells.Item(13, 4).Value = '''650.41'
$ws.Cells.Item(13, 5).Value = '  +8.17%  '
$ws.Cells.Item(14, 4).Value = '3.881.62'
$ws.Cells.Item(14, 5).Value = '  +1.21%  '
$ws.Cells.Item(15, 4).Value = '''8.52'
$ws.Cells.Item(15, 5).Value = '  -1.12%  '
$ws.Cells.Item(16, 4).Value = '66.591.94'
$ws.Cells.Item(16, 5).Value = '  +1.11%  '
$ws.Cells.Item(17, 5).Value = '  +0.49%  '
$ws.Cells.Item(18, 4).Value = '''17.95'
$ws.Cells.Item(18, 5).Value = '  +0.54%  '
$ws.Cells.Item(19, 4).Value = '3.346.01'
$ws.Cells.Item(19, 5).Value = '  +1.11%  '
$ws.Cells.Item(20, 5).Value = '  +1.17%  '
$ws.Cells.Item(21, 5).Value = '  +0.39%  '
$ws.Cells.Item(22, 4).Value = '''17.73'
$ws.Cells.Item(22, 5).Value = '  -4.51%  '
$ws.Cells.Item(23, 4).Value = '''5.09'
$ws.Cells.Item(23, 5).Value = '  +0.71%  '
$ws.Cells.Item(24, 4).Value = '''100.37'
$ws.Cells.Item(24, 5).Value = '  +0.00%  '
$ws.Cells.Item(25, 5).Value = '  +1.84%  '
$ws.Cells.Item(26, 4).Value = '''2.81'
$ws.Cells.Item(26, 5).Value = '  +2.65%  '
$ws.Cells.Item(27, 4).Value = '''9.67'
$ws.Cells.Item(27, 5).Value = '  +1.86%  '
$ws.Cells.Item(28, 4).Value = '''32.14'
$ws.Cells.Item(28, 5).Value = '  +6.43%  '
$ws.Cells.Item(29, 4).Value = '''8.56'
$ws.Cells.Item(29, 5).Value = '  -1.16%  '
$ws.Cells.Item(30, 4).Value = '''6.98'
$ws.Cells.Item(30, 5).Value = '  +3.54%  '
$ws.Cells.Item(31, 2).Value = 'Bittensor'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(31, 4).Value = '''611.30'
$ws.Cells.Item(31, 5).Value = '  +7.84%  '
$ws.Cells.Item(32, 2).Value = 'dogwifhat'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(32, 4).Value = '''3.96'
$ws.Cells.Item(32, 5).Value = '  +1.91%  '
$ws.Cells.Item(33, 4).Value = '''11.13'
$ws.Cells.Item(33, 5).Value = '  +1.31%  '
$ws.Cells.Item(34, 4).Value = '3.887.97'
$ws.Cells.Item(34, 5).Value = '  +5.30%  '
$ws.Cells.Item(35, 5).Value = '  +1.35%  '
$ws.Cells.Item(36, 5).Value = '  +0.12%  '
$ws.Cells.Item(37, 4).Value = '''56.08'
$ws.Cells.Item(37, 5).Value = '  -1.39%  '
$ws.Cells.Item(38, 4).Value = '''2.76'
$ws.Cells.Item(38, 5).Value = '  +4.48%  '
$ws.Cells.Item(39, 5).Value = '  +1.52%  '
$ws.Cells.Item(40, 4).Value = '''33.41'
$ws.Cells.Item(40, 5).Value = '  -1.07%  '
$ws.Cells.Item(41, 2).Value = 'Stacks'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(41, 4).Value = '''3.21'
$ws.Cells.Item(41, 5).Value = '  -1.16%  '
$ws.Cells.Item(42, 2).Value = 'PEPE'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(42, 4).Value = '0.0₃0699'
$ws.Cells.Item(42, 5).Value = '  -1.13%  '
$ws.Cells.Item(43, 4).Value = '''0.343'
$ws.Cells.Item(43, 5).Value = '  +1.68%  '
$ws.Cells.Item(44, 4).Value = '''3.38'
$ws.Cells.Item(44, 5).Value = '  -0.92%  '
$ws.Cells.Item(45, 4).Value = '''0.0419'
$ws.Cells.Item(45, 5).Value = '  +0.51%  '
$ws.Cells.Item(46, 5).Value = '  -0.20%  '
$ws.Cells.Item(47, 2).Value = 'ThetaToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(47, 4).Value = '''2.56'
$ws.Cells.Item(47, 5).Value = '  -0.14%  '
$ws.Cells.Item(48, 2).Value = 'Mantle'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(48, 4).Value = '''1.39'
$ws.Cells.Item(48, 5).Value = '  +11.19%  '
$ws.Cells.Item(49, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(49, 4).Value = '''1.00'
$ws.Cells.Item(49, 5).Value = '  +0.30%  '
$ws.Cells.Item(50, 4).Value = '''2.87'
$ws.Cells.Item(50, 5).Value = '  -16.75%  '
$ws.Cells.Item(51, 4).Value = '''129.26'
$ws.Cells.Item(51, 5).Value = '  +5.22%  '
